# data creation system change and added some coin data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up existing "details" text in column E (rows 2-7): drop macrons
# from "Brahmi" and replace a comma before "holding"/"facing" with "and".
$ws.Range("E2").Value = 'King standing and holding an Indian long bow in his left hand and an arrow in his right Garuda standard at left;Brahmi legend under arm: skanda'
$ws.Range("E3").Value = 'Lakshmi seated facing and holding long-stemmed lotus and diadem;Brahmi legend at right: sri skandaguptah;circular Brahmi legend around'
$ws.Range("E4").Value = 'King standing and holding an Indian long bow in his left hand and an arrow in his right Garuda standard at left;Brahmi legend under arm: skanda;circular Brahmi legend around'
$ws.Range("E5").Value = 'Lakshmi seated facing and holding long-stemmed lotus and diadem;Brahmi legend at right: kramadityah'
$ws.Range("E6").Value = 'King standing at left and facing right and holding an Indian long bow in his left hand and an arrow in his right hand, Lakshmi standing at right, facing left and offering an indistinct object to the king Garuda standard between the two figures,'
$ws.Range("E7").Value = 'Lakshmi seated facing and holding long-stemmed lotus and diadem;Brahmi legend at right: sri skandaguptah;circular Brahmi legend around'

# --- Append newly catalogued coin rows 8-12 ---
$newRows = @(
    @(6,  "Kumara-4829v-625.13-obverse.jpg", 1, 0, "King standing and holding an Indian long bow in his left hand and an arrow in his right Garuda standard at left;Brahmi legend in right field: Kumara", 0),
    @(7,  "Kumara-4829v-625.13-reverse.jpg", 1, 0, "Lakshmi seated facing and holding long-stemmed lotus and scattering coins with her right hand;Brahmi legend at right: Sri Mahendra;circular Brahmi legend around", 0),
    @(8,  "Kumara-4830-280.60-obverse.jpg", 1, 0, "King standing and holding an Indian long bow in his left hand and an arrow in his right Garuda standard at left;Brahmi legend under arm: Ku;circular Brahmi legend around", 0),
    @(9,  "Kumara-4830-280.60-reverse.jpg", 1, 0, "Lakshmi seated facing and holding long-stemmed lotus and diadem;Brahmi legend at right: Sri Mahendra;", 0),
    @(10, "Kumara-4834-486.22-obverse.jpg", 1, 0, "King standing and holding an Indian long bow in his left hand and an arrow in his right Garuda standard at left, Brāhmī legend under arm: Ku", 0)
)

$r = 8
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}
